$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the taxon data between row 5 and row 6 (columns A, B, E, F, G, H, Z, AB)

$colsToSwap = @("A", "B", "E", "F", "G", "H", "Z", "AB")

foreach ($col in $colsToSwap) {
    $cell5 = $ws.Range("$col`5")
    $cell6 = $ws.Range("$col`6")
    $v5 = $cell5.Value2
    $v6 = $cell6.Value2
    $cell5.Value2 = $v6
    $cell6.Value2 = $v5
}
